# The tracker sheet's last data row is duplicated into a new row immediately
# below it (same blank entry columns, same Username/Date stamp), matching a
# resubmission of the row by the same user on the same date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1
if ($lastRow -lt 1) { $lastRow = 1 }

$newRow = $lastRow + 1

$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
$srcRange.Copy($dstRange)
